# Applies the "Updated cryptos list" data refresh (Sun Jan 14 2024) to Sheet1.
# D-column "Price" and E-column "Volume(1h)" cells are plain text (not numbers),
# so plain-numeric-looking new values are pushed through Set-TextValue to stop
# Excel's COM layer from auto-coercing them into real numbers (which would drop
# formatting like trailing zeros). Non-numeric-looking text is just assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '42.909.72'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.544.79'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  -0.07%  '
Set-TextValue "D5" '306.54'
$ws.Range("E5").Value = '  +2.01%  '
Set-TextValue "D6" '99.44'
$ws.Range("E6").Value = '  +7.21%  '
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("E8").Value = '  +0.05%  '
Set-TextValue "D9" '0.548'
$ws.Range("E9").Value = '  -0.52%  '
Set-TextValue "D10" '37.27'
$ws.Range("E10").Value = '  +2.30%  '
$ws.Range("E11").Value = '  +1.23%  '
Set-TextValue "D12" '7.78'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").Value = '2.932.72'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").Value = '2.529.42'
$ws.Range("E15").Value = '  -0.80%  '
Set-TextValue "D16" '15.25'
$ws.Range("E16").Value = '  +7.04%  '
Set-TextValue "D17" '0.879'
$ws.Range("D18").Value = '42.947.97'
$ws.Range("E18").Value = '  -0.37%  '
Set-TextValue "D19" '13.13'
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("E21").Value = '  -0.19%  '
Set-TextValue "D22" '71.78'
$ws.Range("E22").Value = '  -0.51%  '
Set-TextValue "D23" '254.69'
$ws.Range("E23").Value = '  -0.11%  '
Set-TextValue "D24" '2.95'
$ws.Range("E24").Value = '  +0.91%  '
Set-TextValue "D25" '2.06'
$ws.Range("E25").Value = '  -3.45%  '
Set-TextValue "D26" '27.77'
$ws.Range("E26").Value = '  -4.31%  '
Set-TextValue "D27" '0.999'
$ws.Range("E27").Value = '  -0.13%  '
Set-TextValue "D28" '2.30'
$ws.Range("E28").Value = '  +8.79%  '
Set-TextValue "D29" '10.22'
$ws.Range("E29").Value = '  -0.44%  '
Set-TextValue "D30" '38.88'
$ws.Range("E30").Value = '  +5.03%  '
Set-TextValue "D31" '6.19'
$ws.Range("E31").Value = '  +1.45%  '
Set-TextValue "D32" '157.83'
$ws.Range("E32").Value = '  +3.08%  '
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("E35").Value = '  -2.71%  '
Set-TextValue "D36" '18.83'
$ws.Range("E36").Value = '  +6.72%  '
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("E38").Value = '  +0.53%  '
Set-TextValue "D39" '24.32'
$ws.Range("E39").Value = '  +4.64%  '
Set-TextValue "D41" '2.11'
$ws.Range("E41").Value = '  -4.59%  '
Set-TextValue "D42" '3.46'
$ws.Range("E42").Value = '  +0.55%  '
Set-TextValue "D43" '3.90'
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("D45").Value = '2.072.55'
$ws.Range("E45").Value = '  -1.56%  '
$ws.Range("E46").Value = '  -0.01%  '
Set-TextValue "D47" '86.34'
$ws.Range("E47").Value = '  +1.48%  '
Set-TextValue "D48" '9.03'
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("D49").Value = '2.788.19'
$ws.Range("E49").Value = '  -0.82%  '
Set-TextValue "D50" '0.193'
$ws.Range("E50").Value = '  +1.57%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue "D51" '73.68'
$ws.Range("E51").Value = '  +1.14%  '
